# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets.
# F3: 87 -> 88, F7: 123 -> 125, F9: 320 -> 322

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 88
    $ws.Range("F7").Value = 125
    $ws.Range("F9").Value = 322
}
